# Scheduled market-data refresh: update currentAveragePrice / Leve price /
# Leve profit columns (H:N) for a handful of leve rows across several
# sheets, reflecting newly fetched Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 342.55554
$ws.Range("I9").Value = 126.916664
$ws.Range("J9").Value = 773.8333
$ws.Range("K9").Value = 126.916664
$ws.Range("L9").Value = 773.8333
$ws.Range("M9").Value = 42.083336
$ws.Range("N9").Value = -1111.8333

$ws.Range("H17").Value = 1195.3
$ws.Range("J17").Value = 1195.3
$ws.Range("L17").Value = 3585.9
$ws.Range("N17").Value = -3921.9

$ws.Range("H18").Value = 3422.25
$ws.Range("I18").Value = 1769.7142
$ws.Range("J18").Value = 14990
$ws.Range("K18").Value = 1769.7142
$ws.Range("L18").Value = 14990
$ws.Range("M18").Value = -1485.7142
$ws.Range("N18").Value = -15558

$ws.Range("H53").Value = 3998.5334
$ws.Range("I53").Value = 5043.636
$ws.Range("K53").Value = 5043.636
$ws.Range("M53").Value = -4406.636

$ws.Range("H106").Value = 9165.5
$ws.Range("I106").Value = 8998
$ws.Range("J106").Value = 9199
$ws.Range("K106").Value = 8998
$ws.Range("L106").Value = 9199
$ws.Range("M106").Value = -8367
$ws.Range("N106").Value = -10461

$ws.Range("H107").Value = 2210.7234
$ws.Range("I107").Value = 1950.919
$ws.Range("J107").Value = 3172
$ws.Range("K107").Value = 1950.919
$ws.Range("L107").Value = 3172
$ws.Range("M107").Value = -30.9190000000001
$ws.Range("N107").Value = -7012

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 43631
$ws.Range("J44").Value = 38174.668
$ws.Range("L44").Value = 38174.668
$ws.Range("N44").Value = -39150.668

$ws.Range("H45").Value = 3776.1892
$ws.Range("I45").Value = 4522.241
$ws.Range("J45").Value = 1071.75
$ws.Range("K45").Value = 4522.241
$ws.Range("L45").Value = 1071.75
$ws.Range("M45").Value = -4145.241
$ws.Range("N45").Value = -1825.75

$ws.Range("H55").Value = 57998.1
$ws.Range("J55").Value = 63331.223
$ws.Range("L55").Value = 63331.223
$ws.Range("N55").Value = -63961.223

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 1667.5555
$ws.Range("I122").Value = 1659.4584
$ws.Range("J122").Value = 1732.3334
$ws.Range("K122").Value = 4978.3752
$ws.Range("L122").Value = 5197.0002
$ws.Range("M122").Value = -2528.3752
$ws.Range("N122").Value = -10097.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1934.174
$ws.Range("I31").Value = 997.2245
$ws.Range("J31").Value = 4229.7
$ws.Range("K31").Value = 997.2245
$ws.Range("L31").Value = 4229.7
$ws.Range("M31").Value = -702.2245
$ws.Range("N31").Value = -4819.7

$ws.Range("H34").Value = 1934.174
$ws.Range("I34").Value = 997.2245
$ws.Range("J34").Value = 4229.7
$ws.Range("K34").Value = 997.2245
$ws.Range("L34").Value = 4229.7
$ws.Range("M34").Value = -795.2245
$ws.Range("N34").Value = -4633.7

$ws.Range("H58").Value = 6758.476
$ws.Range("I58").Value = 1930.9286
$ws.Range("J58").Value = 16413.572
$ws.Range("K58").Value = 1930.9286
$ws.Range("L58").Value = 16413.572
$ws.Range("M58").Value = -1727.9286
$ws.Range("N58").Value = -16819.572

$ws.Range("H136").Value = 6758.476
$ws.Range("I136").Value = 1930.9286
$ws.Range("J136").Value = 16413.572
$ws.Range("K136").Value = 5792.7858
$ws.Range("L136").Value = 49240.716
$ws.Range("M136").Value = -3242.7858
$ws.Range("N136").Value = -54340.716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1921.2727
$ws.Range("I14").Value = 1921.2727
$ws.Range("K14").Value = 5763.8181
$ws.Range("M14").Value = -5590.8181

$ws.Range("H113").Value = 29413200
$ws.Range("I113").Value = 665.8
$ws.Range("J113").Value = 41668424
$ws.Range("K113").Value = 1997.4
$ws.Range("L113").Value = 125005272
$ws.Range("M113").Value = 172.6000000000001
$ws.Range("N113").Value = -125009612

$ws.Range("H116").Value = 14982
$ws.Range("I116").Value = 14982
$ws.Range("K116").Value = 44946
$ws.Range("M116").Value = -41504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 33339258
$ws.Range("I113").Value = 41673452
$ws.Range("J113").Value = 2481.1667
$ws.Range("K113").Value = 41673452
$ws.Range("L113").Value = 2481.1667
$ws.Range("M113").Value = -41671282
$ws.Range("N113").Value = -6821.1667

$ws.Range("H122").Value = 61982.824
$ws.Range("I122").Value = 65731.75
$ws.Range("K122").Value = 197195.25
$ws.Range("M122").Value = -194745.25

$ws.Range("H132").Value = 940201.25
$ws.Range("I132").Value = 970427.1
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 2911281.3
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -2908751.3
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1363.1786
$ws.Range("I16").Value = 1265.5186
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1265.5186
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1095.5186
$ws.Range("N16").Value = -4340

$ws.Range("H31").Value = 9244.079
$ws.Range("I31").Value = 6627.393
$ws.Range("J31").Value = 16570.8
$ws.Range("K31").Value = 6627.393
$ws.Range("L31").Value = 16570.8
$ws.Range("M31").Value = -6379.393
$ws.Range("N31").Value = -17066.8

$ws.Range("H46").Value = 3138.25
$ws.Range("I46").Value = 849.5
$ws.Range("J46").Value = 3596
$ws.Range("K46").Value = 849.5
$ws.Range("L46").Value = 3596
$ws.Range("M46").Value = -661.5
$ws.Range("N46").Value = -3972

$ws.Range("H55").Value = 1357.2858
$ws.Range("I55").Value = 357.3125
$ws.Range("J55").Value = 2690.5833
$ws.Range("K55").Value = 357.3125
$ws.Range("L55").Value = 2690.5833
$ws.Range("M55").Value = -184.3125
$ws.Range("N55").Value = -3036.5833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 787.84375
$ws.Range("I113").Value = 940.375
$ws.Range("J113").Value = 330.25
$ws.Range("K113").Value = 2821.125
$ws.Range("L113").Value = 990.75
$ws.Range("M113").Value = -651.125
$ws.Range("N113").Value = -5330.75

$ws.Range("H132").Value = 1779.7142
$ws.Range("I132").Value = 1821.375
$ws.Range("K132").Value = 5464.125
$ws.Range("M132").Value = -2934.125
